$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.376.18'
$ws.Range("E2").Value = '  +0.70%  '
$ws.Range("D3").Value = '1.911.94'
$ws.Range("E3").Value = '  +0.96%  '
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.733'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +11.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '256.35'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +4.22%  '
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.29'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.27%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.370'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +7.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '53.18'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.81%  '
$ws.Range("E11").Value = '  +6.29%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0990'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.20%  '
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '13.07'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +6.71%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '2.186.61'
$ws.Range("E14").Value = '  +0.82%  '
$ws.Range("E15").Value = '  +5.91%  '
$ws.Range("E16").Value = '  +4.38%  '
$ws.Range("D17").Value = '1.916.66'
$ws.Range("E17").Value = '  +0.42%  '
$ws.Range("D18").Value = '35.367.93'
$ws.Range("E18").Value = '  +0.72%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '75.30'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +4.67%  '
$ws.Range("E20").Value = '  +4.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '245.55'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '13.15'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +6.19%  '
$ws.Range("E23").Value = '  +7.61%  '
$ws.Range("E24").Value = '  -0.23%  '
$ws.Range("E25").Value = '  +7.46%  '
$ws.Range("E26").Value = '  +3.21%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '166.35'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.40%  '
$ws.Range("E28").Value = '  +4.23%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.86'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +3.22%  '
$ws.Range("E30").Value = '  +5.19%  '
$ws.Range("D31").Value = '4.127.39'
$ws.Range("E31").Value = '  +0.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.39'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +6.66%  '
$ws.Range("B33").Value = 'TrustWalletToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.65'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +24.87%  '
$ws.Range("B34").Value = 'WEMIXToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.99'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +14.22%  '
$ws.Range("E35").Value = '  +5.85%  '
$ws.Range("E36").Value = '  +5.20%  '
$ws.Range("E37").Value = '  -0.24%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.920'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.29%  '
$ws.Range("E39").Value = '  +1.67%  '
$ws.Range("E40").Value = '  +6.52%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '17.12'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +6.76%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '97.67'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +9.22%  '
$ws.Range("E43").Value = '  +3.28%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0644'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.23%  '
$ws.Range("E45").Value = '  +4.87%  '
$ws.Range("D46").Value = '1.343.84'
$ws.Range("E46").Value = '  +0.89%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.75'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +4.49%  '
$ws.Range("E49").Value = '  -0.36%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '45.15'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -8.36%  '
